$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 110.8604276666667
$ws.Range("H2").Value = 332.581283
$ws.Range("I2").Value = 0.2509786052589675
$ws.Range("J2").Value = 0.2509786052589675
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 141.0704656666667
$ws.Range("N2").Value = 423.211397
$ws.Range("O2").Value = 0.05525296614535039
$ws.Range("P2").Value = 0.05525296614535039
$ws.Range("Q2").Value = 15639.13215494248
$ws.Range("R2").Value = 140752.1893944824
$ws.Range("S2").Value = 0.01386731237958099
$ws.Range("T2").Value = 0.01386731237958099

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 110.8604276666667
$ws.Range("H3").Value = 332.581283
$ws.Range("I3").Value = 0.2509786052589675
$ws.Range("J3").Value = 0.2509786052589675
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 3.309554666666667
$ws.Range("N3").Value = 9.928663999999999
$ws.Range("O3").Value = 0.001296250856544298
$ws.Range("P3").Value = 0.001296250856544298
$ws.Range("Q3").Value = 366.8986457328791
$ws.Range("R3").Value = 3302.087811595912
$ws.Range("S3").Value = 0.00032533123204123
$ws.Range("T3").Value = 0.00032533123204123

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 110.8604276666667
$ws.Range("H4").Value = 332.581283
$ws.Range("I4").Value = 0.2509786052589675
$ws.Range("J4").Value = 0.2509786052589675
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 2408.188354666666
$ws.Range("N4").Value = 7224.565063999999
$ws.Range("O4").Value = 0.9432133721485603
$ws.Range("P4").Value = 0.9432133721485604
$ws.Range("Q4").Value = 266972.790900233
$ws.Range("R4").Value = 2402755.118102097
$ws.Range("S4").Value = 0.2367263766034531
$ws.Range("T4").Value = 0.2367263766034531

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 110.8604276666667
$ws.Range("H5").Value = 332.581283
$ws.Range("I5").Value = 0.2509786052589675
$ws.Range("J5").Value = 0.2509786052589675
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 0.6061513333333334
$ws.Range("N5").Value = 1.818454
$ws.Range("O5").Value = 0.0002374108495449545
$ws.Range("P5").Value = 0.0002374108495449545
$ws.Range("Q5").Value = 67.19819604405356
$ws.Range("R5").Value = 604.783764396482
$ws.Range("S5").Value = 0.00005958504389213925
$ws.Range("T5").Value = 0.00005958504389213924

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 184.841802
$ws.Range("H6").Value = 554.525406
$ws.Range("I6").Value = 0.4184661617850055
$ws.Range("J6").Value = 0.4184661617850055
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 141.0704656666667
$ws.Range("N6").Value = 423.211397
$ws.Range("O6").Value = 0.05525296614535039
$ws.Range("P6").Value = 0.05525296614535039
$ws.Range("Q6").Value = 26075.7190828058
$ws.Range("R6").Value = 234681.4717452522
$ws.Range("S6").Value = 0.02312149667008163
$ws.Range("T6").Value = 0.02312149667008163

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 184.841802
$ws.Range("H7").Value = 554.525406
$ws.Range("I7").Value = 0.4184661617850055
$ws.Range("J7").Value = 0.4184661617850055
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 3.309554666666667
$ws.Range("N7").Value = 9.928663999999999
$ws.Range("O7").Value = 0.001296250856544298
$ws.Range("P7").Value = 0.001296250856544298
$ws.Range("Q7").Value = 611.7440484041761
$ws.Range("R7").Value = 5505.696435637584
$ws.Range("S7").Value = 0.0005424371206486183
$ws.Range("T7").Value = 0.0005424371206486183

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 184.841802
$ws.Range("H8").Value = 554.525406
$ws.Range("I8").Value = 0.4184661617850055
$ws.Range("J8").Value = 0.4184661617850055
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 2408.188354666666
$ws.Range("N8").Value = 7224.565063999999
$ws.Range("O8").Value = 0.9432133721485603
$ws.Range("P8").Value = 0.9432133721485604
$ws.Range("Q8").Value = 445133.8750320017
$ws.Range("R8").Value = 4006204.875288016
$ws.Range("S8").Value = 0.3947028795873
$ws.Range("T8").Value = 0.3947028795873

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 184.841802
$ws.Range("H9").Value = 554.525406
$ws.Range("I9").Value = 0.4184661617850055
$ws.Range("J9").Value = 0.4184661617850055
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 0.6061513333333334
$ws.Range("N9").Value = 1.818454
$ws.Range("O9").Value = 0.0002374108495449545
$ws.Range("P9").Value = 0.0002374108495449545
$ws.Range("Q9").Value = 112.042104738036
$ws.Range("R9").Value = 1008.378942642324
$ws.Range("S9").Value = 0.00009934840697519451
$ws.Range("T9").Value = 0.0000993484069751945

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 95.23175666666667
$ws.Range("H10").Value = 285.69527
$ws.Range("I10").Value = 0.2155966197102082
$ws.Range("J10").Value = 0.2155966197102082
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 141.0704656666667
$ws.Range("N10").Value = 423.211397
$ws.Range("O10").Value = 0.05525296614535039
$ws.Range("P10").Value = 0.05525296614535039
$ws.Range("Q10").Value = 13434.38825922136
$ws.Range("R10").Value = 120909.4943329922
$ws.Range("S10").Value = 0.01191235272990012
$ws.Range("T10").Value = 0.01191235272990012

$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 95.23175666666667
$ws.Range("H11").Value = 285.69527
$ws.Range("I11").Value = 0.2155966197102082
$ws.Range("J11").Value = 0.2155966197102082
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 3.309554666666667
$ws.Range("N11").Value = 9.928663999999999
$ws.Range("O11").Value = 0.001296250856544298
$ws.Range("P11").Value = 0.001296250856544298
$ws.Range("Q11").Value = 315.1747046910311
$ws.Range("R11").Value = 2836.57234221928
$ws.Range("S11").Value = 0.0002794673029674128
$ws.Range("T11").Value = 0.0002794673029674128

$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 95.23175666666667
$ws.Range("H12").Value = 285.69527
$ws.Range("I12").Value = 0.2155966197102082
$ws.Range("J12").Value = 0.2155966197102082
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 2408.188354666666
$ws.Range("N12").Value = 7224.565063999999
$ws.Range("O12").Value = 0.9432133721485603
$ws.Range("P12").Value = 0.9432133721485604
$ws.Range("Q12").Value = 229336.0073991164
$ws.Range("R12").Value = 2064024.066592047
$ws.Range("S12").Value = 0.2033536147006962
$ws.Range("T12").Value = 0.2033536147006963

$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 95.23175666666667
$ws.Range("H13").Value = 285.69527
$ws.Range("I13").Value = 0.2155966197102082
$ws.Range("J13").Value = 0.2155966197102082
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 0.6061513333333334
$ws.Range("N13").Value = 1.818454
$ws.Range("O13").Value = 0.0002374108495449545
$ws.Range("P13").Value = 0.0002374108495449545
$ws.Range("Q13").Value = 57.72485627917556
$ws.Range("R13").Value = 519.52370651258
$ws.Range("S13").Value = 0.00005118497664442101
$ws.Range("T13").Value = 0.00005118497664442101

$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 50.778675
$ws.Range("H14").Value = 152.336025
$ws.Range("I14").Value = 0.1149586132458188
$ws.Range("J14").Value = 0.1149586132458188
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 141.0704656666667
$ws.Range("N14").Value = 423.211397
$ws.Range("O14").Value = 0.05525296614535039
$ws.Range("P14").Value = 0.05525296614535039
$ws.Range("Q14").Value = 7163.371328186326
$ws.Range("R14").Value = 64470.34195367694
$ws.Range("S14").Value = 0.006351804365787654
$ws.Range("T14").Value = 0.006351804365787655

$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 50.778675
$ws.Range("H15").Value = 152.336025
$ws.Range("I15").Value = 0.1149586132458188
$ws.Range("J15").Value = 0.1149586132458188
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 3.309554666666667
$ws.Range("N15").Value = 9.928663999999999
$ws.Range("O15").Value = 0.001296250856544298
$ws.Range("P15").Value = 0.001296250856544298
$ws.Range("Q15").Value = 168.0548008134
$ws.Range("R15").Value = 1512.4932073206
$ws.Range("S15").Value = 0.0001490152008870373
$ws.Range("T15").Value = 0.0001490152008870373

$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 50.778675
$ws.Range("H16").Value = 152.336025
$ws.Range("I16").Value = 0.1149586132458188
$ws.Range("J16").Value = 0.1149586132458188
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 2408.188354666666
$ws.Range("N16").Value = 7224.565063999999
$ws.Range("O16").Value = 0.9432133721485603
$ws.Range("P16").Value = 0.9432133721485604
$ws.Range("Q16").Value = 122284.6138004034
$ws.Range("R16").Value = 1100561.524203631
$ws.Range("S16").Value = 0.1084305012571109
$ws.Range("T16").Value = 0.1084305012571109

$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 50.778675
$ws.Range("H17").Value = 152.336025
$ws.Range("I17").Value = 0.1149586132458188
$ws.Range("J17").Value = 0.1149586132458188
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 0.6061513333333334
$ws.Range("N17").Value = 1.818454
$ws.Range("O17").Value = 0.0002374108495449545
$ws.Range("P17").Value = 0.0002374108495449545
$ws.Range("Q17").Value = 30.77956155615
$ws.Range("R17").Value = 277.01605400535
$ws.Range("S17").Value = 0.00002729242203319969
$ws.Range("T17").Value = 0.00002729242203319969
